$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "25.815.53"
$ws.Range("E2").Value = "  +0.10%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.633.34"
$ws.Range("E3").Value = "  +0.31%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.03%  "

# Row 5 - BNB
$ws.Range("D5").Value = "214.72"
$ws.Range("E5").Value = "  -0.22%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.29%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.08%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -0.23%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -0.37%  "

# Row 10 - Solana
$ws.Range("D10").Value = "19.88"
$ws.Range("E10").Value = "  +2.25%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.0781"
$ws.Range("E11").Value = "  +0.32%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.663.29"
$ws.Range("E12").Value = "  +2.45%  "

# Row 13 - Polkadot
$ws.Range("E13").Value = "  -0.24%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "1.858.56"
$ws.Range("E14").Value = "  +0.37%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "0.557"
$ws.Range("E15").Value = "  +0.15%  "

# Row 16 - ShibaInu
$ws.Range("D16").Value = "0.0₃0767"
$ws.Range("E16").Value = "  +1.68%  "

# Row 17 - Litecoin
$ws.Range("E17").Value = "  -0.10%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "25.824.39"
$ws.Range("E18").Value = "  +0.11%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  -0.09%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "193.93"
$ws.Range("E20").Value = "  +0.17%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "4.38"
$ws.Range("E21").Value = "  +1.92%  "

# Row 22 - Avalanche
$ws.Range("E22").Value = "  +1.09%  "

# Row 23 - Chainlink
$ws.Range("D23").Value = "6.18"
$ws.Range("E23").Value = "  +3.01%  "

# Row 24 - BinanceUSD
$ws.Range("E24").Value = "  -0.07%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  -2.68%  "

# Row 26 - Monero
$ws.Range("D26").Value = "139.64"
$ws.Range("E26").Value = "  -0.67%  "

# Row 28 - Cosmos
$ws.Range("E28").Value = "  +1.55%  "

# Row 29 - EthereumClassic
$ws.Range("E29").Value = "  +0.96%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +0.19%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  +1.48%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  +0.99%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  +2.16%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  +1.18%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  +0.50%  "

# Row 36 - ARBITRUM
$ws.Range("E36").Value = "  +0.97%  "

# Row 38 - ImmutableX
$ws.Range("D38").Value = "0.550"
$ws.Range("E38").Value = "  +0.88%  "

# Row 39 - Maker
$ws.Range("D39").Value = "1.120.60"
$ws.Range("E39").Value = "  -1.11%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  +0.42%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  +0.78%  "

# Row 42 - FraxShare
$ws.Range("E42").Value = "  -0.71%  "

# Row 43 - Quant
$ws.Range("D43").Value = "99.62"
$ws.Range("E43").Value = "  +2.39%  "

# Row 44 - TrustWalletToken
$ws.Range("E44").Value = "  +0.58%  "

# Row 45 - BabyDogeCoin
$ws.Range("D45").Value = "0.0₆0109"
$ws.Range("E45").Value = "  -3.42%  "

# Row 46 - Aave
$ws.Range("D46").Value = "55.43"
$ws.Range("E46").Value = "  +0.68%  "

# Row 47 - Mantle
$ws.Range("E47").Value = "  -4.81%  "

# Row 48 - EnergySwap
$ws.Range("D48").Value = "7.70"
$ws.Range("E48").Value = "  +1.49%  "

# Row 49 - Cronos
$ws.Range("E49").Value = "  -0.46%  "

# Row 50 - was SynthetixNetwork, now Frax
$ws.Range("B50").Value = "Frax"
$ws.Range("C50").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.22%  "

# Row 51 - was Frax, now SynthetixNetwork
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").Value = "2.32"
$ws.Range("E51").Value = "  +5.58%  "
